# Fruta / hortaliza, semanal
#
# The source data (Terminal Hortofruticola Agro Chillan - Cilantro) was
# re-synced: for each data row (2-21) the "record payload" - Fecha (D),
# Volumen (J), Precio minimo (K), Precio maximo (L), Precio promedio
# ponderado (M), Unidad de comercializacion (N), Origen (O), Precio $/Kg
# (P) and Kg o Unidades (Q) - now lines up with a different row than
# before, while the market/category columns (A, B, C, E-I, R), which are
# identical on every row, stay put. Row 16 already matched and needs no
# change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Cells.Item(2, 4).Value = 44267   # Fecha
$ws.Cells.Item(2, 10).Value = 150   # Volumen
$ws.Cells.Item(2, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(2, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(2, 13).Value = 1913   # Precio promedio ponderado
$ws.Cells.Item(2, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(2, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(2, 16).Value = 1913   # Precio $/Kg
$ws.Cells.Item(2, 17).Value = 1   # Kg o Unidades

# Row 3
$ws.Cells.Item(3, 4).Value = 44525   # Fecha
$ws.Cells.Item(3, 10).Value = 60   # Volumen
$ws.Cells.Item(3, 11).Value = 2000   # Precio minimo
$ws.Cells.Item(3, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(3, 13).Value = 2000   # Precio promedio ponderado
$ws.Cells.Item(3, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(3, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(3, 16).Value = 2000   # Precio $/Kg
$ws.Cells.Item(3, 17).Value = 1   # Kg o Unidades

# Row 4
$ws.Cells.Item(4, 4).Value = 44273   # Fecha
$ws.Cells.Item(4, 10).Value = 140   # Volumen
$ws.Cells.Item(4, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(4, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(4, 13).Value = 1914   # Precio promedio ponderado
$ws.Cells.Item(4, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(4, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(4, 16).Value = 1914   # Precio $/Kg
$ws.Cells.Item(4, 17).Value = 1   # Kg o Unidades

# Row 5
$ws.Cells.Item(5, 4).Value = 44266   # Fecha
$ws.Cells.Item(5, 10).Value = 150   # Volumen
$ws.Cells.Item(5, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(5, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(5, 13).Value = 1913   # Precio promedio ponderado
$ws.Cells.Item(5, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(5, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(5, 16).Value = 1913   # Precio $/Kg
$ws.Cells.Item(5, 17).Value = 1   # Kg o Unidades

# Row 6
$ws.Cells.Item(6, 4).Value = 44270   # Fecha
$ws.Cells.Item(6, 10).Value = 260   # Volumen
$ws.Cells.Item(6, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(6, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(6, 13).Value = 1908   # Precio promedio ponderado
$ws.Cells.Item(6, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(6, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(6, 16).Value = 1908   # Precio $/Kg
$ws.Cells.Item(6, 17).Value = 1   # Kg o Unidades

# Row 7
$ws.Cells.Item(7, 4).Value = 44211   # Fecha
$ws.Cells.Item(7, 10).Value = 120   # Volumen
$ws.Cells.Item(7, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(7, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(7, 13).Value = 1883   # Precio promedio ponderado
$ws.Cells.Item(7, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(7, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(7, 16).Value = 1883   # Precio $/Kg
$ws.Cells.Item(7, 17).Value = 1   # Kg o Unidades

# Row 8
$ws.Cells.Item(8, 4).Value = 44533   # Fecha
$ws.Cells.Item(8, 10).Value = 100   # Volumen
$ws.Cells.Item(8, 11).Value = 2000   # Precio minimo
$ws.Cells.Item(8, 12).Value = 2200   # Precio maximo
$ws.Cells.Item(8, 13).Value = 2100   # Precio promedio ponderado
$ws.Cells.Item(8, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(8, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(8, 16).Value = 2100   # Precio $/Kg
$ws.Cells.Item(8, 17).Value = 1   # Kg o Unidades

# Row 9
$ws.Cells.Item(9, 4).Value = 44260   # Fecha
$ws.Cells.Item(9, 10).Value = 220   # Volumen
$ws.Cells.Item(9, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(9, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(9, 13).Value = 1909   # Precio promedio ponderado
$ws.Cells.Item(9, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(9, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(9, 16).Value = 1909   # Precio $/Kg
$ws.Cells.Item(9, 17).Value = 1   # Kg o Unidades

# Row 10
$ws.Cells.Item(10, 4).Value = 44263   # Fecha
$ws.Cells.Item(10, 10).Value = 140   # Volumen
$ws.Cells.Item(10, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(10, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(10, 13).Value = 1914   # Precio promedio ponderado
$ws.Cells.Item(10, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(10, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(10, 16).Value = 1914   # Precio $/Kg
$ws.Cells.Item(10, 17).Value = 1   # Kg o Unidades

# Row 11
$ws.Cells.Item(11, 4).Value = 44532   # Fecha
$ws.Cells.Item(11, 10).Value = 100   # Volumen
$ws.Cells.Item(11, 11).Value = 2000   # Precio minimo
$ws.Cells.Item(11, 12).Value = 2200   # Precio maximo
$ws.Cells.Item(11, 13).Value = 2100   # Precio promedio ponderado
$ws.Cells.Item(11, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(11, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(11, 16).Value = 2100   # Precio $/Kg
$ws.Cells.Item(11, 17).Value = 1   # Kg o Unidades

# Row 12
$ws.Cells.Item(12, 4).Value = 44265   # Fecha
$ws.Cells.Item(12, 10).Value = 220   # Volumen
$ws.Cells.Item(12, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(12, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(12, 13).Value = 1909   # Precio promedio ponderado
$ws.Cells.Item(12, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(12, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(12, 16).Value = 1909   # Precio $/Kg
$ws.Cells.Item(12, 17).Value = 1   # Kg o Unidades

# Row 13
$ws.Cells.Item(13, 4).Value = 44539   # Fecha
$ws.Cells.Item(13, 10).Value = 60   # Volumen
$ws.Cells.Item(13, 11).Value = 2000   # Precio minimo
$ws.Cells.Item(13, 12).Value = 2200   # Precio maximo
$ws.Cells.Item(13, 13).Value = 2100   # Precio promedio ponderado
$ws.Cells.Item(13, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(13, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(13, 16).Value = 2100   # Precio $/Kg
$ws.Cells.Item(13, 17).Value = 1   # Kg o Unidades

# Row 14
$ws.Cells.Item(14, 4).Value = 44166   # Fecha
$ws.Cells.Item(14, 10).Value = 240   # Volumen
$ws.Cells.Item(14, 11).Value = 600   # Precio minimo
$ws.Cells.Item(14, 12).Value = 700   # Precio maximo
$ws.Cells.Item(14, 13).Value = 633   # Precio promedio ponderado
$ws.Cells.Item(14, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(14, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(14, 16).Value = 633   # Precio $/Kg
$ws.Cells.Item(14, 17).Value = 1   # Kg o Unidades

# Row 15
$ws.Cells.Item(15, 4).Value = 44271   # Fecha
$ws.Cells.Item(15, 10).Value = 200   # Volumen
$ws.Cells.Item(15, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(15, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(15, 13).Value = 1920   # Precio promedio ponderado
$ws.Cells.Item(15, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(15, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(15, 16).Value = 1920   # Precio $/Kg
$ws.Cells.Item(15, 17).Value = 1   # Kg o Unidades

# Row 16: payload unchanged, nothing to write.

# Row 17
$ws.Cells.Item(17, 4).Value = 44208   # Fecha
$ws.Cells.Item(17, 10).Value = 130   # Volumen
$ws.Cells.Item(17, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(17, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(17, 13).Value = 1908   # Precio promedio ponderado
$ws.Cells.Item(17, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(17, 15).Value = "Provincia de Cautín"   # Origen
$ws.Cells.Item(17, 16).Value = 1908   # Precio $/Kg
$ws.Cells.Item(17, 17).Value = 1   # Kg o Unidades

# Row 18
$ws.Cells.Item(18, 4).Value = 44524   # Fecha
$ws.Cells.Item(18, 10).Value = 80   # Volumen
$ws.Cells.Item(18, 11).Value = 2000   # Precio minimo
$ws.Cells.Item(18, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(18, 13).Value = 2000   # Precio promedio ponderado
$ws.Cells.Item(18, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(18, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(18, 16).Value = 2000   # Precio $/Kg
$ws.Cells.Item(18, 17).Value = 1   # Kg o Unidades

# Row 19
$ws.Cells.Item(19, 4).Value = 44264   # Fecha
$ws.Cells.Item(19, 10).Value = 130   # Volumen
$ws.Cells.Item(19, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(19, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(19, 13).Value = 1908   # Precio promedio ponderado
$ws.Cells.Item(19, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(19, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(19, 16).Value = 1908   # Precio $/Kg
$ws.Cells.Item(19, 17).Value = 1   # Kg o Unidades

# Row 20
$ws.Cells.Item(20, 4).Value = 44272   # Fecha
$ws.Cells.Item(20, 10).Value = 150   # Volumen
$ws.Cells.Item(20, 11).Value = 1800   # Precio minimo
$ws.Cells.Item(20, 12).Value = 2000   # Precio maximo
$ws.Cells.Item(20, 13).Value = 1893   # Precio promedio ponderado
$ws.Cells.Item(20, 14).Value = "`$/atado 0,5 a 1 kilo"   # Unidad de comercializacion
$ws.Cells.Item(20, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(20, 16).Value = 1893   # Precio $/Kg
$ws.Cells.Item(20, 17).Value = 1   # Kg o Unidades

# Row 21
$ws.Cells.Item(21, 4).Value = 44160   # Fecha
$ws.Cells.Item(21, 10).Value = 190   # Volumen
$ws.Cells.Item(21, 11).Value = 1300   # Precio minimo
$ws.Cells.Item(21, 12).Value = 1500   # Precio maximo
$ws.Cells.Item(21, 13).Value = 1395   # Precio promedio ponderado
$ws.Cells.Item(21, 14).Value = "`$/atado 1 a 1,5 kilos"   # Unidad de comercializacion
$ws.Cells.Item(21, 15).Value = "Provincia de Diguillín"   # Origen
$ws.Cells.Item(21, 16).Value = 930   # Precio $/Kg
$ws.Cells.Item(21, 17).Value = 1.5   # Kg o Unidades
